$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete client rows (function: excluir cliente). Rows being removed (in original sheet):
#   Row 2 -> "Correção"
#   Row 5 -> "addteste"
#   Row 6 -> "nova adição"
#   Row 8 -> "NOVO"
# Delete from bottom to top so row indices of not-yet-deleted rows stay valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()
